# Generate Report for Handoff
# b.md has now been handed off: update its status from
# "Handed back: in sync with en-US" to "Ready for handoff" on all three
# sheets, refresh the handback file/date for b.md's rows on the zh-cn and
# de-de sheets, flip the "Content Duplicate" flag to False, and add a
# warning about the handback file being out of date.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

# --- Overview sheet: row 3 is b.md ---
$ws1.Range("E3").Value = "Ready for handoff"
$ws1.Range("F3").Value = "Ready for handoff"
$ws1.Range("G3").Value = "2016-08-15 14:36:37"

# --- zh-cn sheet: row 3 is b.md ---
$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("C3").Value = "Ready for handoff"
# "False" is stored as plain text (not a native boolean) in this workbook,
# so force text entry with a quote prefix, then strip the prefix style back
# off so the cell ends up as an ordinary, unstyled text cell.
$ws2.Range("F3").Value = "'False"
$ws2.Range("F3").Style = "Normal"
$ws2.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-15 14:36:33"
$ws2.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/29e3d8d038653015cf5d0610901190bb61706f5a/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9d5f23c07928552cb584c228ab3308b14678f0ba/e2e/b.md."

# --- de-de sheet: row 3 is b.md ---
$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("F3").Value = "'False"
$ws3.Range("F3").Style = "Normal"
$ws3.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-15 14:36:37"
$ws3.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/29e3d8d038653015cf5d0610901190bb61706f5a/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9d5f23c07928552cb584c228ab3308b14678f0ba/e2e/b.md."

# The new, longer error text in column P ("Error Detail") needs a wider
# column on both localized sheets.
$ws2.Range("P1").ColumnWidth = 39.166666666666664
$ws3.Range("P1").ColumnWidth = 39.166666666666664
